$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert 13 new rows before row 2. This shifts the existing data
# (previously rows 2-21) down to rows 15-34.
$ws.Rows("2:14").Insert()

# Excel's row insert copies the formatting of the row above (the bold
# header row) into the newly inserted rows. The source data rows are
# unstyled, so strip that inherited formatting back off.
$ws.Rows("2:14").ClearFormats()

# Step 2: Remove the trailing 3 rows that fall beyond the new extent (row 31),
# i.e. what used to be rows 19-21 (timestamps 1700/1800/1900), which are no
# longer part of the dataset.
$ws.Rows("32:34").Delete()

# Step 3: Fill the 13 newly inserted rows (2-14) with the new sensor data.
$newData = @(
    @(0, "falling", -3.616065740585328, 4.469793319702149, 0.258730050176382, 0.0594066455960273, 0.1411098688840866, 0.0326812900602817),
    @(100, "falling", -3.719920873641968, 4.70874035358429, 0.04408367723226582, 0.0723875313997268, 0.0080939643085002, 0.0835358202457428),
    @(200, "falling", -3.721616864204406, 4.585918724536896, 0.2726323418319225, 0.0029016099870204, 0.0274889357388019, 0.0937678143382072),
    @(300, "falling", -3.805010795593262, 4.544945240020752, 0.311984956264496, 0.0433714315295219, 0.0073303831741213, -0.0807869285345077),
    @(400, "falling", -3.887511849403381, 4.449418604373932, 0.4409204423427582, 0.0352774672210216, 0.0056505035609006, 0.0056505035609006),
    @(500, "falling", -3.778247833251953, 4.410304188728333, 0.5136718302965164, 0.0557414554059505, 0.0113010071218013, -0.1820378452539444),
    @(600, "falling", -3.584390580654144, 4.580866992473602, 0.3528684750199313, -0.0462730415165424, -0.0134390350431203, -0.0310014113783836),
    @(700, "falling", -3.46419882774353, 4.518833875656127, 0.403674334287644, 0.0276416521519422, 0.001527163083665, -0.0335975885391235),
    @(800, "falling", -3.414171874523162, 4.371547281742096, 0.5634630396962166, -0.0464257597923278, -0.0105374250560998, -0.06368270516395561),
    @(900, "falling", -3.423850417137146, 4.383899688720703, 0.5505108982324599, -0.0219911485910415, -0.0183259565383195, 0.0233655963093042),
    @(1000, "falling", -3.575843572616578, 4.333066165447235, 0.460273951292038, -0.0684169083833694, -0.0335975885391235, 0.0587957799434661),
    @(1100, "falling", -3.668661117553711, 4.494052410125732, 0.1578152179718018, 0.0198531206697225, -0.0583376325666904, 0.0022907445672899),
    @(1200, "falling", -3.737768590450287, 4.342036247253418, 0.1853629685938359, -0.0027488935738801, -0.0503963828086853, 0.0137444678694009)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = 2 + $i
    $rowVals = $newData[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# Step 4: Fix up column A (timestamp) for all the rows that were shifted
# down by the insert (old rows 2-18, now at rows 15-31); the timestamp
# sequence must keep counting up by 100 (1300, 1400, ... 2900) rather than
# keeping the old values that came along with the shift (0, 100, ... 1600).
for ($r = 15; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
}
